$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.250.25'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.892.07'
$ws.Range('E3').Value = '  -0.40%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9994'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.23%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '322.51'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.81%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.15%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4707'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +2.42%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.4028'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.61%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '47.35'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.93%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.08014'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.15%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.9965'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.27%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '22.69'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('D13').Value = '1.934.56'
$ws.Range('E13').Value = '  +3.20%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.891'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.82%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.044'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.21%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '89.38'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.9997'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.29%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.06632'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.94%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.00001021'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.93%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.48'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '29.233.16'
$ws.Range('E22').Value = '  +0.43%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.501'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +2.19%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.177'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').Value = '2.084.62'
$ws.Range('E26').Value = '  -0.71%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '155.07'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.71%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '19.67'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.45%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.024'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +6.10%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.084'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -2.03%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '118.97'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.47%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.033'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.84%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.09411'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.388'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.539'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.40%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.363'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.19%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.06066'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.02232'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.42%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.169'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.82%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '8.022'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.87%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.5821'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.492'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +8.51%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.1829'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.06%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '10.05'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.07726'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.49%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.250'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.87%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '12.09'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.52%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.5476'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.84%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.903'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.29%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '113.58'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.30%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '44.02'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.07%  '
